# Append the missing invoice row (row 6) to the 2025 invoices sheet,
# correcting the fiscal-year data set (commit: "base.html update correct fiscal year").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 6
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")

# Force every new cell to be treated as text so values such as leading-zero
# AFMs, dates written as plain strings, and numeric-looking totals are not
# auto-converted to numbers/dates by Excel.
foreach ($col in $cols) {
    $ws.Range($col + $newRow).NumberFormat = "@"
}

$ws.Range("A$newRow").Value = "400011186892779"
$ws.Range("B$newRow").Value = "094439854"
$ws.Range("C$newRow").Value = "ΤΡΑΚΑΔΑΣ Α.Ε."
$ws.Range("D$newRow").Value = "8Μ0ΤΔΑ"
$ws.Range("E$newRow").Value = "8970"
$ws.Range("F$newRow").Value = "04/10/2025"
$ws.Range("G$newRow").Value = "Τιμολόγιο Πώλησης"
$ws.Range("H$newRow").Value = ""
$ws.Range("I$newRow").Value = "48,39"
$ws.Range("J$newRow").Value = "11,61"
$ws.Range("K$newRow").Value = "60,00"
$ws.Range("L$newRow").Value = "8970"

# Match the plain (unstyled) look of the rest of the data rows instead of
# keeping the temporary text-number-format style used above.
$ws.Range("A$newRow" + ":L$newRow").Style = $ws.Range("A5:L5").Style
